$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed assignment for D-column values that would otherwise be
# auto-coerced to numbers by Excel (e.g. "0.999", "130.18"), then restore the
# default (Normal) cell style so no stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "56.272.92"
$ws.Cells.Item(2, 5).Value = "  -2.32%  "
$ws.Cells.Item(3, 4).Value = "2.374.19"
$ws.Cells.Item(3, 5).Value = "  -1.86%  "
$ws.Cells.Item(4, 5).Value = "  -0.24%  "
Set-TextValue $ws.Cells.Item(5, 4) "501.19"
$ws.Cells.Item(5, 5).Value = "  -1.44%  "
Set-TextValue $ws.Cells.Item(6, 4) "130.18"
$ws.Cells.Item(6, 5).Value = "  -2.47%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.999"
$ws.Cells.Item(7, 5).Value = "  +0.46%  "
$ws.Cells.Item(8, 5).Value = "  -2.81%  "
$ws.Cells.Item(9, 4).Value = "2.382.61"
$ws.Cells.Item(9, 5).Value = "  -3.13%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.0984"
$ws.Cells.Item(10, 5).Value = "  -0.34%  "
$ws.Cells.Item(11, 5).Value = "  +0.61%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.324"
$ws.Cells.Item(12, 5).Value = "  +0.36%  "
Set-TextValue $ws.Cells.Item(13, 4) "4.66"
$ws.Cells.Item(13, 5).Value = "  -0.20%  "
$ws.Cells.Item(14, 4).Value = "2.796.86"
$ws.Cells.Item(14, 5).Value = "  -2.46%  "
$ws.Cells.Item(15, 4).Value = "56.232.52"
$ws.Cells.Item(15, 5).Value = "  -2.19%  "
Set-TextValue $ws.Cells.Item(16, 4) "21.59"
$ws.Cells.Item(16, 5).Value = "  -1.94%  "
$ws.Cells.Item(17, 5).Value = "  -1.68%  "
$ws.Cells.Item(18, 4).Value = "2.407.77"
$ws.Cells.Item(18, 5).Value = "  -3.49%  "
$ws.Cells.Item(19, 5).Value = "  -2.88%  "
$ws.Cells.Item(20, 5).Value = "  -3.22%  "
Set-TextValue $ws.Cells.Item(21, 4) "307.48"
$ws.Cells.Item(21, 5).Value = "  -2.83%  "
Set-TextValue $ws.Cells.Item(22, 4) "6.25"
$ws.Cells.Item(22, 5).Value = "  -1.94%  "
$ws.Cells.Item(23, 5).Value = "  +0.47%  "
Set-TextValue $ws.Cells.Item(24, 4) "64.72"
$ws.Cells.Item(24, 5).Value = "  -1.27%  "
$ws.Cells.Item(25, 5).Value = "  -0.10%  "
Set-TextValue $ws.Cells.Item(26, 4) "0.370"
$ws.Cells.Item(26, 5).Value = "  -4.02%  "
$ws.Cells.Item(27, 5).Value = "  -4.72%  "
Set-TextValue $ws.Cells.Item(28, 4) "7.29"
$ws.Cells.Item(28, 5).Value = "  -4.59%  "
Set-TextValue $ws.Cells.Item(29, 4) "172.18"
$ws.Cells.Item(29, 5).Value = "  -1.15%  "
$ws.Cells.Item(30, 5).Value = "  -3.36%  "
$ws.Cells.Item(31, 5).Value = "  -3.90%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.999"
$ws.Cells.Item(32, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(33, 4) "5.75"
$ws.Cells.Item(33, 5).Value = "  -7.90%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.998"
$ws.Cells.Item(34, 5).Value = "  +0.55%  "
$ws.Cells.Item(35, 5).Value = "  -5.29%  "
Set-TextValue $ws.Cells.Item(36, 4) "17.59"
$ws.Cells.Item(36, 5).Value = "  -2.55%  "
$ws.Cells.Item(37, 5).Value = "  -6.46%  "
Set-TextValue $ws.Cells.Item(38, 4) "3.77"
$ws.Cells.Item(38, 5).Value = "  -3.19%  "
$ws.Cells.Item(39, 5).Value = "  -1.31%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.792"
$ws.Cells.Item(40, 5).Value = "  -4.12%  "
$ws.Cells.Item(41, 5).Value = "  -4.16%  "
Set-TextValue $ws.Cells.Item(42, 4) "131.10"
$ws.Cells.Item(42, 5).Value = "  -3.33%  "
$ws.Cells.Item(43, 5).Value = "  -2.62%  "
Set-TextValue $ws.Cells.Item(44, 4) "4.76"
$ws.Cells.Item(44, 5).Value = "  -5.46%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.566"
$ws.Cells.Item(45, 5).Value = "  -1.51%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.0903"
$ws.Cells.Item(46, 5).Value = "  -1.92%  "
Set-TextValue $ws.Cells.Item(47, 4) "241.88"
$ws.Cells.Item(47, 5).Value = "  -7.42%  "
$ws.Cells.Item(48, 5).Value = "  -3.28%  "
$ws.Cells.Item(49, 5).Value = "  -3.05%  "
Set-TextValue $ws.Cells.Item(50, 4) "16.98"
$ws.Cells.Item(50, 5).Value = "  -1.97%  "
$ws.Cells.Item(51, 5).Value = "  -3.14%  "
